# "Generate Report for Handback"
#
# The localization-status report is refreshed: the zh-de "handback" run is
# recorded (zh-cn handback timestamp + per-row "Latest Target File" /
# "Latest Handback File" links), a new de-de handback run is recorded with
# its own timestamp, the Overview sheet's status text flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and a few
# columns that now hold long file names / hyperlinks are widened so the
# text is readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: status text + wider zh-cn / de-de columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value() = "Handed back: in sync with en-US"
$overview.Range("F2").Value() = "Handed back: in sync with en-US"
$overview.Range("E3").Value() = "Handed back: in sync with en-US"
$overview.Range("F3").Value() = "Handed back: in sync with en-US"

# Target stored column width is 29.9777047293527 chars; this runtime's
# ColumnWidth setter only lands on 1/6-character boundaries (stored =
# round(ColumnWidth*6)/6 + 5/6), so the closest reachable stored width is
# 30 - use the middle of the input range that rounds to it.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# 2. zh-cn sheet: handback already completed previously -> fill in the
#    Latest Target File / Latest Handback File / Latest Handback DateTime
#    columns and add matching hyperlinks for the new "Latest Target File"
#    cells (column I), mirroring column A's hyperlinks.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

$zhcn.Range("I2").Value() = "2545a184-5ca5-4e30-949d-554e46b672f9.md"
$zhcn.Range("J2").Value() = "2545a184-5ca5-4e30-949d-554e46b672f9.04d56327738f6a8aa55d368838971d14f0dce0a4.zh-cn.xlf"
$zhcn.Range("K2").Value() = "2016-09-01 17:11:57"

$zhcn.Range("I3").Value() = "9bd0504e-e9ab-4178-8393-021487888d92.md"
$zhcn.Range("J3").Value() = "9bd0504e-e9ab-4178-8393-021487888d92.a4e5b1c10e2638e108b9babcc9f18a0b426c3d62.zh-cn.xlf"
$zhcn.Range("K3").Value() = "2016-09-01 17:11:57"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/2545a184-5ca5-4e30-949d-554e46b672f9.md", "", "", "2545a184-5ca5-4e30-949d-554e46b672f9.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/2545a184-5ca5-4e30-949d-554e46b672f9.md", "", "", "2545a184-5ca5-4e30-949d-554e46b672f9.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/9bd0504e-e9ab-4178-8393-021487888d92.md", "", "", "9bd0504e-e9ab-4178-8393-021487888d92.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/9bd0504e-e9ab-4178-8393-021487888d92.md", "", "", "9bd0504e-e9ab-4178-8393-021487888d92.md")

# ---------------------------------------------------------------------
# 3. de-de sheet: this is the language that just got handed back - fill
#    in the same columns with its own (later) handback timestamp and add
#    matching hyperlinks.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

$dede.Range("I2").Value() = "2545a184-5ca5-4e30-949d-554e46b672f9.md"
$dede.Range("J2").Value() = "2545a184-5ca5-4e30-949d-554e46b672f9.04d56327738f6a8aa55d368838971d14f0dce0a4.de-de.xlf"
$dede.Range("K2").Value() = "2016-09-01 17:12:15"

$dede.Range("I3").Value() = "9bd0504e-e9ab-4178-8393-021487888d92.md"
$dede.Range("J3").Value() = "9bd0504e-e9ab-4178-8393-021487888d92.a4e5b1c10e2638e108b9babcc9f18a0b426c3d62.de-de.xlf"
$dede.Range("K3").Value() = "2016-09-01 17:12:15"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/2545a184-5ca5-4e30-949d-554e46b672f9.md", "", "", "2545a184-5ca5-4e30-949d-554e46b672f9.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/2545a184-5ca5-4e30-949d-554e46b672f9.md", "", "", "2545a184-5ca5-4e30-949d-554e46b672f9.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/9bd0504e-e9ab-4178-8393-021487888d92.md", "", "", "9bd0504e-e9ab-4178-8393-021487888d92.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c855c3c3a2e69742b660b0aa9dba33d91dd0aa8a/e2e/9bd0504e-e9ab-4178-8393-021487888d92.md", "", "", "9bd0504e-e9ab-4178-8393-021487888d92.md")
